# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for two more worker rows: shift the previously-last data row
#    (old row 21, the bottom-bordered row) down to row 23, and shift the
#    footer ("firma" block, old rows 26/27) down to rows 28/29.
#    We do this by copying formats explicitly (instead of Rows.Insert, which
#    would mint brand-new style records) so the workbook's style table is
#    reused exactly like the real edit did.
# ---------------------------------------------------------------------------

# Footer block: move old row 27 -> row 29, old row 26 -> row 28.
$ws.Range("B27:J27").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B26:J26").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Clear the now-stale footer rows before we overwrite row 23 below (row 26/27
# sit below the new row 23 and must end up blank).
$ws.Range("B26:J27").ClearContents()
$ws.Range("B26:J27").ClearFormats()

# Bottom-bordered last data row: old row 21 (style ids 21-26) becomes row 23.
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New middle rows 21 & 22 get the same "plain" row style as row 20.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Header summary numbers
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 416188   # VALOR MORA total
$ws.Range("C13").Value2 = 8        # Cant. Trabajadores

# ---------------------------------------------------------------------------
# 3) Worker data rows 16-23 (period 2506 -> 2507, new workers, new amounts)
# ---------------------------------------------------------------------------
function Set-WorkerRow {
    param($row, $doc, $name, $mora, $salario)
    $ws.Range("B$row").Value2 = "CC"
    $ws.Range("C$row").Value2 = $doc
    $ws.Range("D$row").Value2 = $name
    $ws.Range("E$row").Value2 = "2507"
    $ws.Range("F$row").Value2 = $mora
    $ws.Range("G$row").Value2 = $salario
}

Set-WorkerRow 16 "72285009"   "CARLOS ANDRES ACOSTA TERAN"      63160 900000
Set-WorkerRow 17 "1047455386" "MARIA GRACIELA PALACIO VEGA"     56940 908526
Set-WorkerRow 18 "1048936436" "ALFONSO LUIS RODRIGUEZ PADILLA"  56940 1423500
Set-WorkerRow 19 "72284047"   "ELKIN ALBERTO NUÑEZ SORACA"      56940 900000
Set-WorkerRow 20 "1235046610" "DAIRIS MILETH PINEDA ROJAS"      56940 1300000
Set-WorkerRow 21 "1050969759" "EDINSON MANUEL GOMEZ OLIVARES"   11388 1423500
Set-WorkerRow 22 "85290132"   "OSCAR LUIS TORRES RODRIGUEZ"     56940 1423500
Set-WorkerRow 23 "1235046104" "SEBASTIAN TAPIA PALACIO"         56940 1423500

# ---------------------------------------------------------------------------
# 4) Dimension / used range bookkeeping is handled automatically on save.
# ---------------------------------------------------------------------------
